# ajuste: corrigindo as categorias
# Adds a "Total" column (X) with the row sum for each existing category row,
# then adds two new category rows: "Outros" and "Total" (grand total),
# each with their own age-bracket breakdown and row total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: X1 = "Total" ---
$ws.Cells.Item(1, 24).Value = "Total"

# --- New column X values (row totals) for existing rows 2-6 ---
$ws.Cells.Item(2, 24).Value = 2238
$ws.Cells.Item(3, 24).Value = 286
$ws.Cells.Item(4, 24).Value = 915
$ws.Cells.Item(5, 24).Value = 239
$ws.Cells.Item(6, 24).Value = 1486

# --- New row 7: "Outros" ---
$row7 = @{
    1  = "Outros"
    2  = 294
    3  = 10
    4  = 15
    5  = 85
    6  = 139
    7  = 171
    8  = 211
    9  = 191
    10 = 244
    11 = 325
    12 = 435
    13 = 559
    14 = 705
    15 = 844
    16 = 921
    17 = 1000
    18 = 1094
    19 = 1094
    20 = 694
    21 = 276
    22 = 71
    23 = 2
    24 = 9380
}
foreach ($col in $row7.Keys) {
    $ws.Cells.Item(7, $col).Value = $row7[$col]
}

# --- New row 8: "Total" (grand total across all categories) ---
$row8 = @{
    1  = "Total"
    2  = 312
    3  = 12
    4  = 18
    5  = 92
    6  = 148
    7  = 192
    8  = 242
    9  = 256
    10 = 326
    11 = 464
    12 = 652
    13 = 878
    14 = 1148
    15 = 1400
    16 = 1552
    17 = 1628
    18 = 1802
    19 = 1740
    20 = 1118
    21 = 452
    22 = 110
    23 = 2
    24 = 14544
}
foreach ($col in $row8.Keys) {
    $ws.Cells.Item(8, $col).Value = $row8[$col]
}
